$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "Descriere" column
$ws.Range("B1").Value = "Descriere"
$ws.Range("B1").Font.Bold = $true

# Fill in the description values for each "Tema" row (rows 2-31)
for ($i = 1; $i -le 30; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = "descriere $i"
}

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 11.5546875

# Update the active selection like in the edited workbook
$ws.Range("C23").Select()
